$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 (Row R1 / CONDITIONS): abc/def/ghi rule text, now extended with VAR:{var2} tokens
$c2Text = @"
(
(
( ( {
	"message": "TOKEN: abc, TOKEN: def, TOKEN: ghi, REGEX:[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}-[0-9a-z]{4}, VAR: var1, REGEX:Bar[0-9]{2}-sector[0-9]{2}, VAR: var2 "
} AND {
	"message": "TOKEN:Messaging, VAR:{var2}"
} ) OR ( {
	"message": "TOKEN:Messaging, VAR:{var2}"
} AND {
	"message": "TOKEN:updates, VAR:{var2}"
}) )
) AND
{
	"message": "TOKEN:DistributedMessaging, VAR:{var2}"
}
)
OR
{
	"message": "TOKEN:DistributedMessaging Watchdog - subscriber watchdog looking for updates, VAR:{var2}"
}
"@

# C3 (Row R2 / CONDITIONS): DistributedMessaging/Watchdog rule text (unchanged content)
$c3Text = @"
(
(
( ( {
	"message": "TOKEN:DistributedMessaging"
} AND {
	"message": "TOKEN:Watchdog"
} ) OR ( {
	"message": "TOKEN:subscriber"
} AND {
	"message": "TOKEN:updates"
}) )
) AND
{
	"message": "TOKEN:DistributedMessaging"
}
)
OR
{
	"message": "TOKEN:DistributedMessaging Watchdog - subscriber watchdog looking for updates"
}
"@

$ws.Range("C2").Value = $c2Text
$ws.Range("C3").Value = $c3Text

# Preserve the original manually-set row heights (writing the new, longer
# text triggers an autofit recalculation we don't want to surface as a diff)
$ws.Rows.Item(2).RowHeight = 138.75
$ws.Rows.Item(3).RowHeight = 153

[void]$ws.Range("B2").Select()
